$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 2 (shifts the old row 2 down to row 4)
$ws.Rows.Item(2).Resize(2).Insert()

# ---------------------------------------------------------------------------
# Row 2 (new observation: Campanula cervicaria, id 173360)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 173360
$ws.Range("B2").Value = 104643
$ws.Range("C2").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 245
$ws.Range("F2").Value = "Skogsklocka"
$ws.Range("G2").Value = "Campanula cervicaria"
$ws.Range("H2").Value = "L."

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "260"
$ws.Range("J2").Value = "plantor/tuvor"
$ws.Range("K2").Value = "blomning"

$ws.Range("P2").Value = "Valdemarsviks ishockeyhall, Ög"
$ws.Range("Q2").Value = 593615.080207533
$ws.Range("R2").Value = 6453736.111156315
$ws.Range("S2").Value = 10
$ws.Range("T2").Value = "Östergötland"
$ws.Range("U2").Value = "Valdemarsvik"
$ws.Range("V2").Value = "Östergötland"
$ws.Range("W2").Value = "Valdemarsvik"

$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2011-07-12"
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value = "00:00"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2011-07-12"
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "00:00"

$ws.Range("AC2").Value = "Ett nytt stort bestånd (ca 190 ex) utanför nordvästhörnet av hockeyhallen."

$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false

$ws.Range("AI2").Value = "skogsslänt"

$ws.Range("AW2").Value = "Stefan Kasselstrand"
$ws.Range("AX2").Value = "Stefan Kasselstrand"

# ---------------------------------------------------------------------------
# Row 3 (new observation: Campanula cervicaria, id 167634)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 167634
$ws.Range("B3").Value = 104643
$ws.Range("C3").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 245
$ws.Range("F3").Value = "Skogsklocka"
$ws.Range("G3").Value = "Campanula cervicaria"
$ws.Range("H3").Value = "L."

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "30"
$ws.Range("J3").Value = "plantor/tuvor"
$ws.Range("K3").Value = "överblommad"

$ws.Range("P3").Value = "Valdemarsviks ishall SV, 50 m, Ög"
$ws.Range("Q3").Value = 593635.2000545097
$ws.Range("R3").Value = 6453756.08671589
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Östergötland"
$ws.Range("U3").Value = "Valdemarsvik"
$ws.Range("V3").Value = "Östergötland"
$ws.Range("W3").Value = "Valdemarsvik"
$ws.Range("X3").Value = "E-Val-0147"

$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2010-07-31"
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = "00:00"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2010-07-31"
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value = "00:00"

$ws.Range("AC3").Value = "76P, , identisk med Obs Id 4545496"

$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false

$ws.Range("AW3").Value = "E-län Floraväktarna"
$ws.Range("AX3").Value = "Stefan Kasselstrand"
$ws.Range("AY3").Value = "Floraväkteri Sverige"
